$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add "PD" and "N2" column headers, matching existing header style
$ws.Range("Q1").Value = "PD"
$ws.Range("R1").Value = "N2"

# Copy header formatting (bold, centered, bordered) from P1 onto the new header cells
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-78: PD and N2 values
$ws.Range("Q2").Value = "PUE-F"
$ws.Range("R2").Value = "Fuera de Poligono OVL"
$ws.Range("Q3").Value = "BLO-?"
$ws.Range("R3").Value = "Fuera de Poligono OVL"
$ws.Range("Q4").Value = "BLO-M"
$ws.Range("R4").Value = "Fuera de Poligono OVL"
$ws.Range("Q5").Value = "BLO-F"
$ws.Range("R5").Value = "Fuera de Poligono OVL"
$ws.Range("Q6").Value = "PCH-C"
$ws.Range("R6").Value = "Fuera de Poligono OVL"
$ws.Range("Q7").Value = "BLO-P"
$ws.Range("R7").Value = "Fuera de Poligono OVL"
$ws.Range("Q8").Value = "VCR-K"
$ws.Range("R8").Value = "Fuera de Poligono OVL"
$ws.Range("Q9").Value = "BLO-B"
$ws.Range("R9").Value = "Fuera de Poligono OVL"
$ws.Range("Q10").Value = "ATH-F"
$ws.Range("R10").Value = "Fuera de Poligono OVL"
$ws.Range("Q11").Value = "PPT-P"
$ws.Range("R11").Value = "Fuera de Poligono OVL"
$ws.Range("Q12").Value = "PUE-E"
$ws.Range("R12").Value = "Fuera de Poligono OVL"
$ws.Range("Q13").Value = "BLO-I"
$ws.Range("R13").Value = "Fuera de Poligono OVL"
$ws.Range("Q14").Value = "ATH-A"
$ws.Range("R14").Value = "Fuera de Poligono OVL"
$ws.Range("Q15").Value = "COG-F"
$ws.Range("R15").Value = "Fuera de Poligono OVL"
$ws.Range("Q16").Value = "PUE-I"
$ws.Range("R16").Value = "Fuera de Poligono OVL"
$ws.Range("Q17").Value = "PUE-I"
$ws.Range("R17").Value = "Fuera de Poligono OVL"
$ws.Range("Q18").Value = "COG-O"
$ws.Range("R18").Value = "Fuera de Poligono OVL"
$ws.Range("Q19").Value = "BLO-P"
$ws.Range("R19").Value = "Fuera de Poligono OVL"
$ws.Range("Q20").Value = "ATH-R"
$ws.Range("R20").Value = "Fuera de Poligono OVL"
$ws.Range("Q21").Value = "PUE-K"
$ws.Range("R21").Value = "Fuera de Poligono OVL"
$ws.Range("Q22").Value = "PUE-K"
$ws.Range("R22").Value = "Fuera de Poligono OVL"
$ws.Range("Q23").Value = "BLO-G"
$ws.Range("R23").Value = "Fuera de Poligono OVL"
$ws.Range("Q24").Value = "COG-L"
$ws.Range("R24").Value = "Fuera de Poligono OVL"
$ws.Range("Q25").Value = "COG-A"
$ws.Range("R25").Value = "Fuera de Poligono OVL"
$ws.Range("Q26").Value = "PUE-J"
$ws.Range("R26").Value = "Fuera de Poligono OVL"
$ws.Range("Q27").Value = "ATH-S"
$ws.Range("R27").Value = "Fuera de Poligono OVL"
$ws.Range("Q28").Value = "PUE-G"
$ws.Range("R28").Value = "Fuera de Poligono OVL"
$ws.Range("Q29").Value = "PUE-O"
$ws.Range("R29").Value = "ARATO-25058.PO.1PUE"
$ws.Range("Q30").Value = "AGU-O"
$ws.Range("R30").Value = "Fuera de Poligono OVL"
$ws.Range("Q31").Value = "ATH-P"
$ws.Range("R31").Value = "Fuera de Poligono OVL"
$ws.Range("Q32").Value = "ATH-P"
$ws.Range("R32").Value = "Fuera de Poligono OVL"
$ws.Range("Q33").Value = "COG-B"
$ws.Range("R33").Value = "Fuera de Poligono OVL"
$ws.Range("Q34").Value = "BLO-H"
$ws.Range("R34").Value = "Fuera de Poligono OVL"
$ws.Range("Q35").Value = "ATH-B"
$ws.Range("R35").Value = "Fuera de Poligono OVL"
$ws.Range("Q36").Value = "PCH-K"
$ws.Range("R36").Value = "Fuera de Poligono OVL"
$ws.Range("Q37").Value = "AGU-K"
$ws.Range("R37").Value = "Fuera de Poligono OVL"
$ws.Range("Q38").Value = "PCH-I"
$ws.Range("R38").Value = "Fuera de Poligono OVL"
$ws.Range("Q39").Value = "BLO-L"
$ws.Range("R39").Value = "Fuera de Poligono OVL"
$ws.Range("Q40").Value = "COG-A"
$ws.Range("R40").Value = "Fuera de Poligono OVL"
$ws.Range("Q41").Value = "COG-N"
$ws.Range("R41").Value = "Fuera de Poligono OVL"
$ws.Range("Q42").Value = "ATH-O"
$ws.Range("R42").Value = "Fuera de Poligono OVL"
$ws.Range("Q43").Value = "VCR-N"
$ws.Range("R43").Value = "Fuera de Poligono OVL"
$ws.Range("Q44").Value = "VCR-D"
$ws.Range("R44").Value = "Fuera de Poligono OVL"
$ws.Range("Q45").Value = "AGU-O"
$ws.Range("R45").Value = "Fuera de Poligono OVL"
$ws.Range("Q46").Value = "VCR-I"
$ws.Range("R46").Value = "Fuera de Poligono OVL"
$ws.Range("Q47").Value = "PUE-P"
$ws.Range("R47").Value = "ARATO-25058.PO.1PUE"
$ws.Range("Q48").Value = "PUE-N"
$ws.Range("R48").Value = "ARATO-25058.PO.2PUE"
$ws.Range("Q49").Value = "PUE-A"
$ws.Range("R49").Value = "Fuera de Poligono OVL"
$ws.Range("Q50").Value = "PUE-J"
$ws.Range("R50").Value = "Fuera de Poligono OVL"
$ws.Range("Q51").Value = "VCR-?"
$ws.Range("R51").Value = "Fuera de Poligono OVL"
$ws.Range("Q52").Value = "ATH-H"
$ws.Range("R52").Value = "Fuera de Poligono OVL"
$ws.Range("Q53").Value = "BLO-J"
$ws.Range("R53").Value = "Fuera de Poligono OVL"
$ws.Range("Q54").Value = "VCR-O"
$ws.Range("R54").Value = "Fuera de Poligono OVL"
$ws.Range("Q55").Value = "COG-I"
$ws.Range("R55").Value = "Fuera de Poligono OVL"
$ws.Range("Q56").Value = "BLO-E"
$ws.Range("R56").Value = "Fuera de Poligono OVL"
$ws.Range("Q57").Value = "COG-F"
$ws.Range("R57").Value = "Fuera de Poligono OVL"
$ws.Range("Q58").Value = "COG-I"
$ws.Range("R58").Value = "Fuera de Poligono OVL"
$ws.Range("Q59").Value = "BLO-Q"
$ws.Range("R59").Value = "Fuera de Poligono OVL"
$ws.Range("Q60").Value = "CLI-O"
$ws.Range("R60").Value = "Fuera de Poligono OVL"
$ws.Range("Q61").Value = "ATH-Q"
$ws.Range("R61").Value = "Fuera de Poligono OVL"
$ws.Range("Q62").Value = "CLI-O"
$ws.Range("R62").Value = "Fuera de Poligono OVL"
$ws.Range("Q63").Value = "PUE-A"
$ws.Range("R63").Value = "Fuera de Poligono OVL"
$ws.Range("Q64").Value = "COG-L"
$ws.Range("R64").Value = "Fuera de Poligono OVL"
$ws.Range("Q65").Value = "NRA-F"
$ws.Range("R65").Value = "ARATO-25058.PO.1NRA"
$ws.Range("Q66").Value = "COG-P"
$ws.Range("R66").Value = "Fuera de Poligono OVL"
$ws.Range("Q67").Value = "PUE-K"
$ws.Range("R67").Value = "Fuera de Poligono OVL"
$ws.Range("Q68").Value = "BLO-Q"
$ws.Range("R68").Value = "Fuera de Poligono OVL"
$ws.Range("Q69").Value = "COG-K"
$ws.Range("R69").Value = "Fuera de Poligono OVL"
$ws.Range("Q70").Value = "COG-Q"
$ws.Range("R70").Value = "Fuera de Poligono OVL"
$ws.Range("Q71").Value = "COG-L"
$ws.Range("R71").Value = "Fuera de Poligono OVL"
$ws.Range("Q72").Value = "COG-P"
$ws.Range("R72").Value = "Fuera de Poligono OVL"
$ws.Range("Q73").Value = "NRA-K"
$ws.Range("R73").Value = "Fuera de Poligono OVL"
$ws.Range("Q74").Value = "ATH-A"
$ws.Range("R74").Value = "Fuera de Poligono OVL"
$ws.Range("Q75").Value = "VCR-I"
$ws.Range("R75").Value = "Fuera de Poligono OVL"
$ws.Range("Q76").Value = "VCR-I"
$ws.Range("R76").Value = "Fuera de Poligono OVL"
$ws.Range("Q77").Value = "VCR-I"
$ws.Range("R77").Value = "Fuera de Poligono OVL"
$ws.Range("Q78").Value = "BLO-P"
$ws.Range("R78").Value = "Fuera de Poligono OVL"
